$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N4, N6 value updates (N7 AVERAGE formula will auto-recalculate)
$ws.Range("N4").Value = 72.5
$ws.Range("N6").Value = 72.3

# C22 allocation count correction
$ws.Range("C22").Value = 23

# Row 28 (A28/D28) and Row 43 (A43/D43): fix allocation error - swap/correct track
# geometry text and recompute distance totals. SUM formulas in C29/D29/D49 recalc
# automatically.
$ws.Range("A28").Value = 'LineString (525327.85414873412810266 178420.77501689142081887, 525435.54454665677621961 178546.54447064295527525, 525345.92272875038906932 178653.20884751385892741, 525429.15519386087544262 178703.14832658017985523, 525460.07010947330854833 178656.1814355535316281, 525501.7821286671096459 178615.87244897283380851, 525558.91436985018663108 178563.64555782591924071, 525601.34650842915289104 178515.0539152596029453, 525591.59441469295416027 178508.14618219647672959, 525612.31761388236191124 178515.8665897376195062, 525739.09483245271258056 178641.42479659095988609, 525767.53843918326310813 178664.58601921438821591, 525708.21320228814147413 178737.72672223576228134, 525731.62045715074054897 178716.68428954426781274, 525769.37393091234844178 178666.83504263576469384, 525806.39432751352433115 178612.95387134497286752, 526026.40911028429400176 178722.27400289245997556, 526144.43453311175107956 178765.52555535719147883, 526150.06246284244116396 178767.20975544024258852, 526255.77580733434297144 178814.4714835460181348, 526271.20979337661992759 178817.35269416749360971, 526312.78246261400636286 178730.88154215391841717, 526431.85428454133216292 178507.28957353180157952, 526371.81847308482974768 178458.83962042655912228, 526552.97916730435099453 178570.48516453863703646, 526555.08568700461182743 178631.57423584524076432, 526531.73185393039602786 178797.25923736946424469, 526501.3099303717026487 178977.1575194436009042, 526543.68378180800937116 178800.00690841992036439, 526578.25740370713174343 178795.88277246299549006, 526808.88298217917326838 178842.85059796637506224, 526867.58985495707020164 178857.72481234182487242, 526887.47400313417892903 178869.08718272874830291, 526862.61881791288033128 179084.26207193068694323, 526876.82178089639637619 179026.0299236977880355, 526893.15518832765519619 179011.11681256498559378, 526922.56527758203446865 179021.44788825331488624, 526976.82480131415650249 179026.93843529763398692, 527126.68443828821182251 179083.05505582402111031, 527206.78183046414051205 179136.0226861338596791, 527205.48993704188615084 179152.81730062235146761, 527216.03914421703666449 179145.72578124998835847, 527240.65396095986943692 179168.27212599763879552, 527356.07469829195179045 179281.62439133448060602, 527415.23299735458567739 179328.3718584259331692, 527437.15880050358828157 179344.91963438750826754, 527437.57249490264803171 179363.53588234429480508, 527448.74728261982090771 179356.647314573492622, 527504.23852299549616873 179398.66757797525497153, 527585.06438483926467597 179335.59891305171186104, 527665.89024668303318322 179312.33086191487382166, 527724.67269166023470461 179303.14610488718608394, 527769.37184252834413201 179306.82000769826117903, 527780.39355096162762493 179281.10268802070640959, 527792.02757653000298887 179222.32024304344668053, 527763.24867117661051452 179195.99060623071272857, 527768.86287219019141048 179171.7344985579547938, 527738.35830124607309699 179152.01631094136973843, 527772.52151037787552923 179165.95770141950924881, 527886.39517858054023236 179195.27266738505568355, 527822.45171340426895767 179520.46795567701337859, 527858.1322810254059732 179530.35083859681617469, 527891.88773283653426915 179539.25118037080392241, 527870.7807220573304221 179648.64833412505686283, 527868.61700371222104877 179658.18725142459152266, 527860.88046702521387488 179707.83289611263899133, 527871.11952563317026943 179673.40761047863634303, 527874.78819991135969758 179657.57498131820466369, 527913.67431851325090975 179664.11720327276270837, 527912.70541233289986849 179757.62383855623193085, 527908.07293678901623935 179798.69010824267752469, 527909.8863292447058484 179817.68885372436488979, 527912.46024989418219775 179831.26684675173601136, 527920.7886318500386551 179757.18808093448751606, 528047.62872848019469529 179774.55966792465187609, 528127.87066305975895375 179815.19500658992910758, 528096.22333193966187537 179801.94153121206909418, 528121.1935451285680756 179822.46549322339706123, 528115.72619099263101816 179848.85401823552092537, 528228.20267733931541443 179870.84723329346161336, 528232.45957181649282575 179854.12167384184431285, 528255.63412408879958093 179859.83663536846870556, 528233.04146977188065648 179858.24069136375328526, 528229.61167071235831827 179870.97432596748694777, 528256.6169783161021769 179876.28377280529821292, 528294.42466580052860081 179904.86705374554730952, 528301.72961289319209754 179899.58161686407402158, 528306.33070521731860936 179893.07007217648788355, 528319.42971817124634981 179870.80175015470013022, 528331.87378047755919397 179834.77946453119511716, 528377.07891597924754024 179802.95061608715332113, 528464.50208342960104346 179798.41427575697889552, 528500.39834169461391866 179812.81222550064558163, 528513.98412839835509658 179812.26687349064741284, 528529.45710146625060588 179823.30889789306093007, 528522.1132692665560171 179877.67606056467047893, 528516.61473636329174042 179896.37107243580976501, 528982.89032656128983945 180285.39227534280507825, 528985.75060231669340283 180264.89964467784739099, 528997.18651210982352495 180296.38934114936273545, 529128.14837153907865286 180399.20853238902054727, 529260.66181104525458068 180474.31033623719122261, 529390.66158775170333683 180271.80423245803103782, 529460.4421130547998473 180310.72029464627848938, 529430.42441646789666265 180361.7129266606643796, 529518.94092159741558135 180419.07557802950032055, 529552.56730343436356634 180373.58106142663746141, 529648.71743816568050534 180427.49491339456290007, 529608.480220987345092 180501.82199512678198516, 529569.44679091090802103 180572.25595415246789344, 529667.74284364131744951 180629.72475247565307654, 529730.9295924516627565 180529.48919386579655111, 529787.31384936452377588 180567.2724588074197527, 529763.12702995492145419 180622.33490404754411429, 529706.82561028259806335 180576.71049719117581844, 529752.5037739867111668 180495.92673939117230475, 529740.29997511499095708 180487.27313655486796051, 529754.05698475218378007 180486.60747479822020978, 529800.36167021188884974 180412.51190174778457731, 529824.728592021856457 180415.94968910259194672, 529890.00781108462251723 180459.31649564448161982, 529865.28168949380051345 180503.69017882505431771, 529855.61934717127587646 180535.74483229906763881, 529894.55648332042619586 180548.55958596843993291, 529890.50570647965651006 180564.17520972224883735, 529872.95531774952542037 180606.76134864933555946, 529915.37259634933434427 180613.59366198081988841, 529972.59322050062473863 180618.71789697947679088, 529975.86554363591130823 180653.11958405951736495, 529993.10341755871195346 180649.76780548423994333, 530003.82197457575239241 180653.6264860104129184, 530054.15927445457782596 180611.80309442835277878, 530069.85973392939195037 180508.62456558644771576, 530160.12516808835789561 180516.73132610463653691, 530168.55684169579762965 180494.19840525725157931, 530148.45289235515519977 180461.23712512332713231, 530174.04404229309875518 180496.50880550872534513, 530161.336840910022147 180527.69920890356297605, 530154.59773694002069533 180615.4214020338258706, 530232.82974896696396172 180614.62270036846166477, 530270.34142313164193183 180619.80239770526532084, 530278.97677504969760776 180605.4567557742993813, 530319.42912215983960778 180633.89190510252956301, 530452.98087174398824573 180736.21395450358977541, 530441.01278525067027658 180758.81861569813918322, 530400.55581117747351527 180811.16401171206962317, 530524.57064777321647853 180911.87319211009889841, 530583.04193836136255413 180879.49776674556778744)'
$ws.Range("D28").Value = 10865.1751107169

$ws.Range("A43").Value = 'LineString (528079.47067953646183014 180632.24651449063094333, 528074.47738231881521642 180628.91345400654245168, 528024.52704086224548519 180609.87107211141847074, 528002.13440979458391666 180610.44380836229538545, 527990.64072253345511854 180618.42097470845328644, 527972.56485775776673108 180621.4870105440786574, 527925.40382085286546499 180703.35145196394296363, 527902.95878737303428352 180701.55584928556345403, 527887.69616460683755577 180736.57010151402209885, 527757.73942075879313052 180828.59473878113203682, 527713.74715513840783387 180842.06175886897835881, 527632.04723327199462801 180837.57275217302958481, 527614.09120648819953203 180851.93757360009476542, 527504.55944310687482357 180858.22218297442304902, 527341.38404970872215927 180835.77714949465007521, 527227.36327963136136532 180809.74091065808897838, 527043.76290576683823019 180770.23765173368155956, 526939.61795042059384286 180763.05524102016352117, 526905.50149953132495284 180726.24538611332536675, 526911.78610890568234026 180702.90255129433353432, 526979.12120934505946934 180662.5014910307363607, 527036.5804950533201918 180622.10043076713918708, 527055.43432317627593875 180555.66313166698091663, 527077.87935665610712022 180308.3188627198105678, 527073.39034996007103473 180283.85377622683881782, 527108.40460218861699104 180190.48243695095879957, 527076.98155531683005393 180196.76704632528708316, 527017.72666693024802953 180212.92747043073177338, 526981.81461336265783757 180219.21207980506005697, 526955.77837452606763691 180225.49668917941744439, 526942.31135443819221109 180205.74505971721373498, 526901.91029417456593364 180116.8627271372533869, 526897.6457378133200109 180099.3556010230386164, 526888.66772442148067057 180050.87432870670454577, 526893.15673111740034074 180035.61170594047871418, 526880.58751236868556589 180035.61170594047871418, 526824.02602799970190972 179903.63490907935192809, 526802.70324619382154197 179852.90913341505802236, 526794.6230341411428526 179796.34764904601615854, 526794.26693146175239235 179743.8655021067825146, 526793.01302744331769645 179731.11747791964444332, 526790.71420340961776674 179729.23662189205060713, 526799.15176595037337393 179665.74589264905080199, 526810.75304763286840171 179467.00080663594417274, 526668.71815055573824793 179454.61404235597001389, 526444.82263462664559484 179424.7408409794443287, 526441.88616490841377527 179401.24908323338604532, 526400.95936604985035956 179579.61037377678439952, 526395.80753890855703503 179617.60509894342976622, 526400.31538765714503825 179630.4846667965175584, 526396.45151730126235634 179663.00557562557514757, 526429.93839371926151216 179670.41132714110426605, 526431.87032889726106077 179659.14170526966336183, 526400.31538765714503825 179655.59982411004602909, 526404.17925801302772015 179629.51869920754688792, 526481.13467593537643552 179644.97418063125223853, 526512.36762797913979739 179646.26213741654646583, 526515.90950913866981864 179639.17837509737000801, 526505.28386565984692425 179635.63649393775267527, 526386.30885761685203761 179621.30797470119432546, 526318.94980677636340261 179621.48525741521734744, 526215.34960218577180058 179626.84875574553734623, 526121.02660867478698492 179655.6876124200061895, 526113.32453510130289942 179686.26937513795564882, 526116.04291400965303183 179658.17945975257316604, 526072.0957883259980008 179674.03667005078750663, 526077.53254614258185029 179626.23850758047774434, 526092.24883591232355684 179409.14943758773733862, 526170.58646800776477903 179414.55203290464123711, 526170.91200886678416282 179254.77395809174049646, 526170.02756560104899108 179226.2506627686379943, 526170.69089805032126606 179197.28514581263880245, 526191.0330931645585224 179103.75527045089984313, 526208.50084766466170549 179068.1564290011592675, 526211.74739217478781939 179033.83077586514991708, 526212.5073481819126755 179000.39271154996822588, 526093.19425505725666881 178993.55310748549527489, 526096.23407908587250859 178919.83737479062983766, 526071.91548685671295971 178938.07631896255770698, 526089.39447502139955759 178936.55640694824978709, 526074.95531088532879949 179075.628356259141583, 525857.43331313808448613 179086.440267137310002, 525822.15490559919271618 179081.18603622727096081, 525774.11622299312148243 179039.15218894690042362, 525789.9340731572592631 179008.44695039271027781, 525608.18420976598281413 178985.48173996026162058, 525584.06940557796042413 178987.97679135657381266, 525573.85263970901723951 179052.99129645497305319, 525566.40184021822642535 179179.43401284626452252, 525539.12785214162431657 179262.48159250925527886, 525230.67203460761811584 179148.19059485456091352, 525218.00911157194059342 179136.66408798881457187, 525275.51739498809911311 178989.94413040805375203, 525288.82619905064348131 178973.64458385962643661, 525321.71609151223674417 178916.64969141795882024, 525110.25522415374871343 178791.72511747083626688, 525140.83571881789248437 178743.57710459540248848, 525170.6157625982305035 178737.53536503901705146, 525136.90032081049866974 178734.29459415673045442, 525073.96749363793060184 178845.07109110968303867, 525050.22967286221683025 178837.15848418447421864, 524920.68373622617218643 178769.44129003383568488, 524812.24932546948548406 178896.28909129649400711, 524699.21156708616763353 179025.31069560500327498, 524681.6391244096448645 179051.25509521568892524, 524679.63010795961599797 179077.09913740976480767, 524788.97600055101793259 179149.0372246410115622, 524957.31112467218190432 179248.31178502010880038, 524889.31033757049590349 179366.96949636968201958, 524837.84125263791065663 179392.39398410747526214, 524751.64603811223059893 179345.88577483102562837, 524614.60184811090584844 179570.98550772908492945, 524703.84473427385091782 179436.35361285935505293, 524844.32857223774772137 179505.77558726171264425, 524810.43752961605787277 179598.63431089642108418, 524822.46338344959076494 179627.87900090057519265, 524845.42183167708572 179640.99811417344608344, 524954.20114589820150286 179688.82821464753942564, 524980.43937244394328445 179689.64815922707202844, 525145.27674076042603701 179339.90930881092208438, 525213.61852516769431531 179374.56420798628823832, 525206.26161919755395502 179390.0524310814216733, 525222.23384926444850862 179385.21236136413062923, 525351.1733065313892439 179466.91273819096386433, 525352.72212884086184204 179482.01375570870004594, 525369.8851592862047255 179479.95885097398422658, 525422.48737169418018311 179507.45546200548415072, 525366.53774576925206929 179596.40102116821799427, 525283.64405398862436414 179742.1564941595424898, 525384.085012789350003 179740.24697022794862278, 525418.45644355763215572 179681.43363313551526517, 525419.82130964624229819 179656.85948007801198401)'
$ws.Range("D43").Value = 10624.362508620899

# Restore the active cell selection as recorded after the edit
$ws.Range("N7").Select()
